# Scheduled-runner refresh of Kujata_Profits market-price data.
# All H:N columns in this workbook are plain literal values (no formulas);
# the runner re-pulls market board prices/profit figures per leve row and
# overwrites the stale numbers in place. This script reproduces that
# refresh, row by row, per sheet (sheet name == job class abbreviation).

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 4694.2144
$ws.Range("I28").Value = 4452.5
$ws.Range("K28").Value = 4452.5
$ws.Range("M28").Value = -3967.5

$ws.Range("H39").Value = 642.1111
$ws.Range("I39").Value = 684.875
$ws.Range("J39").Value = 300
$ws.Range("K39").Value = 2054.625
$ws.Range("L39").Value = 900
$ws.Range("M39").Value = -1758.625
$ws.Range("N39").Value = -1492

$ws.Range("H49").Value = 499.2857
$ws.Range("J49").Value = 499.2857
$ws.Range("L49").Value = 1497.8571
$ws.Range("N49").Value = -1769.8571

$ws.Range("H86").Value = 2352.3076
$ws.Range("I86").Value = 2603.25
$ws.Range("J86").Value = 1950.8
$ws.Range("K86").Value = 2603.25
$ws.Range("L86").Value = 1950.8
$ws.Range("M86").Value = -1480.25
$ws.Range("N86").Value = -4196.8

$ws.Range("H89").Value = 2352.3076
$ws.Range("I89").Value = 2603.25
$ws.Range("J89").Value = 1950.8
$ws.Range("K89").Value = 13016.25
$ws.Range("L89").Value = 9754
$ws.Range("M89").Value = -7400.25
$ws.Range("N89").Value = -20986

$ws.Range("H107").Value = 4488.6665
$ws.Range("I107").Value = 3860.25
$ws.Range("K107").Value = 3860.25
$ws.Range("M107").Value = -1940.25

$ws.Range("H109").Value = 28200
$ws.Range("J109").Value = 28200
$ws.Range("L109").Value = 28200
$ws.Range("N109").Value = -30974

$ws.Range("H133").Value = 37275.555
$ws.Range("J133").Value = 37275.555
$ws.Range("L133").Value = 37275.555
$ws.Range("N133").Value = -47395.555

$ws.Range("H137").Value = 2873.5144
$ws.Range("J137").Value = 3002.7407
$ws.Range("L137").Value = 9008.222099999999
$ws.Range("N137").Value = -14108.2221

$ws.Range("H138").Value = 2510.013
$ws.Range("I138").Value = 1788
$ws.Range("J138").Value = 2711.2295
$ws.Range("K138").Value = 5364
$ws.Range("L138").Value = 8133.6885
$ws.Range("M138").Value = -224
$ws.Range("N138").Value = -18413.6885

$ws.Range("H140").Value = 37535.363
$ws.Range("J140").Value = 37535.363
$ws.Range("L140").Value = 37535.363
$ws.Range("N140").Value = -47895.363

$ws.Range("H141").Value = 588.1053000000001
$ws.Range("I141").Value = 588.1053000000001
$ws.Range("K141").Value = 1764.3159
$ws.Range("M141").Value = 3415.6841

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1299.9231
$ws.Range("I45").Value = 1042.7142
$ws.Range("K45").Value = 1042.7142
$ws.Range("M45").Value = -665.7141999999999

$ws.Range("H97").Value = 7916.5
$ws.Range("I97").Value = 838.1818
$ws.Range("J97").Value = 33870.332
$ws.Range("K97").Value = 838.1818
$ws.Range("L97").Value = 33870.332
$ws.Range("M97").Value = -342.1818
$ws.Range("N97").Value = -34862.332

$ws.Range("H110").Value = 517
$ws.Range("I110").Value = 627.75
$ws.Range("J110").Value = 74
$ws.Range("K110").Value = 627.75
$ws.Range("L110").Value = 74
$ws.Range("M110").Value = 1417.25
$ws.Range("N110").Value = -4164

$ws.Range("H124").Value = 17776.428
$ws.Range("J124").Value = 17776.428
$ws.Range("L124").Value = 17776.428
$ws.Range("N124").Value = -27596.428

$ws.Range("H125").Value = 45000
$ws.Range("J125").Value = 45000
$ws.Range("L125").Value = 45000
$ws.Range("N125").Value = -54840

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 168319120
$ws.Range("I105").Value = 201982540
$ws.Range("K105").Value = 201982540
$ws.Range("M105").Value = -201980793

$ws.Range("H107").Value = 1236.5555
$ws.Range("I107").Value = 1254.1428
$ws.Range("J107").Value = 1175
$ws.Range("K107").Value = 1254.1428
$ws.Range("L107").Value = 1175
$ws.Range("M107").Value = 665.8571999999999
$ws.Range("N107").Value = -5015

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 41667980
$ws.Range("I16").Value = 58824812
$ws.Range("K16").Value = 58824812
$ws.Range("M16").Value = -58824525

$ws.Range("H31").Value = 1547.5962
$ws.Range("I31").Value = 1361.2727
$ws.Range("K31").Value = 1361.2727
$ws.Range("M31").Value = -1066.2727

$ws.Range("H34").Value = 1547.5962
$ws.Range("I34").Value = 1361.2727
$ws.Range("K34").Value = 1361.2727
$ws.Range("M34").Value = -1159.2727

$ws.Range("H107").Value = 778.1111
$ws.Range("I107").Value = 469.5238
$ws.Range("J107").Value = 1858.1666
$ws.Range("K107").Value = 469.5238
$ws.Range("L107").Value = 1858.1666
$ws.Range("M107").Value = 1450.4762
$ws.Range("N107").Value = -5698.1666

$ws.Range("H113").Value = 41667980
$ws.Range("I113").Value = 58824812
$ws.Range("K113").Value = 58824812
$ws.Range("M113").Value = -58822642

$ws.Range("H124").Value = 12000
$ws.Range("J124").Value = 12000
$ws.Range("L124").Value = 12000
$ws.Range("N124").Value = -16910

$ws.Range("H141").Value = 356918.12
$ws.Range("J141").Value = 356918.12
$ws.Range("L141").Value = 356918.12
$ws.Range("N141").Value = -367278.12

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1155.9
$ws.Range("I122").Value = 630.7273
$ws.Range("J122").Value = 1459.9474
$ws.Range("K122").Value = 5676.545700000001
$ws.Range("L122").Value = 13139.5266
$ws.Range("M122").Value = -3226.545700000001
$ws.Range("N122").Value = -18039.5266

$ws.Range("H131").Value = 22762256
$ws.Range("J131").Value = 51142.168
$ws.Range("L131").Value = 153426.504
$ws.Range("N131").Value = -163506.504

$ws.Range("H137").Value = 26792086
$ws.Range("I137").Value = 53573520
$ws.Range("J137").Value = 10653.214
$ws.Range("K137").Value = 160720560
$ws.Range("L137").Value = 31959.642
$ws.Range("M137").Value = -160715460
$ws.Range("N137").Value = -42159.642

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1317.1428
$ws.Range("J113").Value = 1400
$ws.Range("L113").Value = 1400
$ws.Range("N113").Value = -5740

$ws.Range("H123").Value = 23585.5
$ws.Range("J123").Value = 23585.5
$ws.Range("L123").Value = 23585.5
$ws.Range("N123").Value = -28485.5

$ws.Range("H138").Value = 37629
$ws.Range("J138").Value = 37629
$ws.Range("L138").Value = 37629
$ws.Range("N138").Value = -47909

$ws.Range("H141").Value = 71779.75
$ws.Range("J141").Value = 71779.75
$ws.Range("L141").Value = 71779.75
$ws.Range("N141").Value = -82139.75

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
# N134 is fully cleared (row no longer has a profit figure), not just zeroed.
$ws.Range("N134").ClearContents()

$ws.Range("H141").Value = 65000
$ws.Range("J141").Value = 65000
$ws.Range("L141").Value = 65000
$ws.Range("N141").Value = -75360

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 250.2
$ws.Range("I113").Value = 176.15384
$ws.Range("J113").Value = 387.7143
$ws.Range("K113").Value = 528.4615200000001
$ws.Range("L113").Value = 1163.1429
$ws.Range("M113").Value = 1641.53848
$ws.Range("N113").Value = -5503.1429
